$d = $word.ActiveDocument
$d.Content.Find.Execute(" tornaria igualmente", $true, $false, $false, $false, $false,
                         $true, 1, $false, " revelaria igualmente", 2)
